# Generate Report for Handback
# Adds "Latest Target File" / "Latest Handback File" hyperlinked entries for
# each localized-file row, and flips the status / handback-timestamp fields
# now that the handback is complete and in sync with en-US.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: every "Ready for handoff" cell (Overview + each language
#    sheet) becomes "Handed back: in sync with en-US" now that handback
#    has completed.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback File"
#    (G) columns for both rows, with hyperlinks matching the existing
#    Source File Name / Latest Handoff File link style, and record the
#    handback datetime.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("F2").Value = "47406448-c21f-42e3-a5d2-dd6660f36f07.md"
$wsZh.Range("G2").Value = "47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-03-23 00:14:26"

$wsZh.Range("F3").Value = "d22991c2-4ab8-454b-8554-17fdf6363a59.md"
$wsZh.Range("G3").Value = "d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-03-23 00:14:26"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/47406448-c21f-42e3-a5d2-dd6660f36f07.md", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33778f20449474906df3811283011d521024858e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.zh-cn.xlf", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/47406448-c21f-42e3-a5d2-dd6660f36f07.md", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33778f20449474906df3811283011d521024858e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.zh-cn.xlf", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/d22991c2-4ab8-454b-8554-17fdf6363a59.md", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33778f20449474906df3811283011d521024858e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.zh-cn.xlf", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/d22991c2-4ab8-454b-8554-17fdf6363a59.md", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33778f20449474906df3811283011d521024858e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.zh-cn.xlf", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.zh-cn.xlf")

$wsZh.Range("F2:G3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, with its own handback datetime.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("F2").Value = "47406448-c21f-42e3-a5d2-dd6660f36f07.md"
$wsDe.Range("G2").Value = "47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-23 00:14:32"

$wsDe.Range("F3").Value = "d22991c2-4ab8-454b-8554-17fdf6363a59.md"
$wsDe.Range("G3").Value = "d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-23 00:14:32"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/47406448-c21f-42e3-a5d2-dd6660f36f07.md", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d05d02b07fbf30c8e1253ec73af32395cb4e3615/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.de-de.xlf", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/47406448-c21f-42e3-a5d2-dd6660f36f07.md", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d05d02b07fbf30c8e1253ec73af32395cb4e3615/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.de-de.xlf", "", "", "47406448-c21f-42e3-a5d2-dd6660f36f07.854a3caf115e97b4ba74804eb4e4988d5b415a55.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/d22991c2-4ab8-454b-8554-17fdf6363a59.md", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d05d02b07fbf30c8e1253ec73af32395cb4e3615/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.de-de.xlf", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/47891ba569871c8d59f2684d0b259e028e2c6c9c/e2e/d22991c2-4ab8-454b-8554-17fdf6363a59.md", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d05d02b07fbf30c8e1253ec73af32395cb4e3615/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.de-de.xlf", "", "", "d22991c2-4ab8-454b-8554-17fdf6363a59.3fdc66c8b68320b253bd646b4afeb758ec19b396.de-de.xlf")

$wsDe.Range("F2:G3").Style = "HyperLink"

Write-Output "Handback report generated."
